# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the 'Mateus Profits' Leve-profit workbook
# (columns H..N: currentAveragePrice*, LevePrice*, LeveProfit*) across all 8 sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: ALC ----
$ws = $wb.Worksheets.Item(1)
# Row 2
$ws.Range("H2").Value = 212.63637
$ws.Range("I2").Value = 183.9
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 183.9
$ws.Range("L2").Value = 500
$ws.Range("M2").Value = -70.90000000000001
$ws.Range("N2").Value = -726
# Row 6
$ws.Range("H6").Value = 1429.8
$ws.Range("I6").Value = 1406.1765
$ws.Range("K6").Value = 4218.529500000001
$ws.Range("M6").Value = -4106.529500000001
# Row 11
$ws.Range("H11").Value = 112.53846
$ws.Range("I11").Value = 112.53846
$ws.Range("K11").Value = 112.53846
$ws.Range("M11").Value = 27.46154
# Row 17
$ws.Range("H17").Value = 4001058
$ws.Range("J17").Value = 4001058
$ws.Range("L17").Value = 12003174
$ws.Range("N17").Value = -12003510
# Row 21
$ws.Range("H21").Value = 5108
$ws.Range("I21").Value = 5108
$ws.Range("K21").Value = 5108
$ws.Range("M21").Value = -4640
# Row 23
$ws.Range("H23").Value = 5108
$ws.Range("I23").Value = 5108
$ws.Range("K23").Value = 5108
$ws.Range("M23").Value = -4874
# Row 62
$ws.Range("H62").Value = 14250.125
$ws.Range("I62").Value = 13418.75
$ws.Range("K62").Value = 13418.75
$ws.Range("M62").Value = -12794.75
# Row 65
$ws.Range("H65").Value = 14250.125
$ws.Range("I65").Value = 13418.75
$ws.Range("K65").Value = 67093.75
$ws.Range("M65").Value = -63973.75
# Row 74
$ws.Range("H74").Value = 8440
$ws.Range("I74").Value = 5017.778
$ws.Range("K74").Value = 5017.778
$ws.Range("M74").Value = -4081.778
# Row 77
$ws.Range("H77").Value = 8440
$ws.Range("I77").Value = 5017.778
$ws.Range("K77").Value = 25088.89
$ws.Range("M77").Value = -20408.89
# Row 86
$ws.Range("H86").Value = 4243.8
$ws.Range("I86").Value = 3998.3333
$ws.Range("K86").Value = 3998.3333
$ws.Range("M86").Value = -2875.3333
# Row 88
$ws.Range("H88").Value = 2715.7646
$ws.Range("I88").Value = 2501.3333
$ws.Range("J88").Value = 2957
$ws.Range("K88").Value = 2501.3333
$ws.Range("L88").Value = 2957
$ws.Range("M88").Value = -2095.3333
$ws.Range("N88").Value = -3769
# Row 89
$ws.Range("H89").Value = 4243.8
$ws.Range("I89").Value = 3998.3333
$ws.Range("K89").Value = 19991.6665
$ws.Range("M89").Value = -14375.6665
# Row 91
$ws.Range("H91").Value = 2715.7646
$ws.Range("I91").Value = 2501.3333
$ws.Range("J91").Value = 2957
$ws.Range("K91").Value = 2501.3333
$ws.Range("L91").Value = 2957
$ws.Range("M91").Value = -1097.3333
$ws.Range("N91").Value = -5765
# Row 98
$ws.Range("H98").Value = 10249.5
$ws.Range("I98").Value = 10249.5
$ws.Range("K98").Value = 10249.5
$ws.Range("M98").Value = -8751.5
# Row 100
$ws.Range("H100").Value = 1627.5385
$ws.Range("I100").Value = 1165.8
$ws.Range("K100").Value = 1165.8
$ws.Range("M100").Value = -624.8
# Row 116
$ws.Range("H116").Value = 3659
$ws.Range("I116").Value = 3230
$ws.Range("J116").Value = 3842.8572
$ws.Range("K116").Value = 3230
$ws.Range("L116").Value = 3842.8572
$ws.Range("M116").Value = 212
$ws.Range("N116").Value = -10726.8572
# Row 122
$ws.Range("H122").Value = 10249.5
$ws.Range("I122").Value = 10249.5
$ws.Range("K122").Value = 30748.5
$ws.Range("M122").Value = -28298.5
# Row 125
$ws.Range("H125").Value = 8586.074000000001
$ws.Range("I125").Value = 9233.846
$ws.Range("J125").Value = 7984.5713
$ws.Range("K125").Value = 83104.614
$ws.Range("L125").Value = 71861.14169999999
$ws.Range("M125").Value = -80644.614
$ws.Range("N125").Value = -76781.14169999999
# Row 129
$ws.Range("H129").Value = 1241
$ws.Range("I129").Value = 1104.1538
$ws.Range("J129").Value = 1834
$ws.Range("K129").Value = 3312.4614
$ws.Range("L129").Value = 5502
$ws.Range("M129").Value = 1687.5386
$ws.Range("N129").Value = -15502
# Row 132
$ws.Range("H132").Value = 2010.8
$ws.Range("I132").Value = 1934.5294
$ws.Range("K132").Value = 5803.5882
$ws.Range("M132").Value = -3273.5882
# Row 135
$ws.Range("H135").Value = 562.15
$ws.Range("I135").Value = 433.8421
$ws.Range("K135").Value = 3904.5789
$ws.Range("M135").Value = -1369.5789
# Row 138
$ws.Range("H138").Value = 16131498
$ws.Range("I138").Value = 2030.8125
$ws.Range("J138").Value = 21741748
$ws.Range("K138").Value = 6092.4375
$ws.Range("L138").Value = 65225244
$ws.Range("M138").Value = -952.4375
$ws.Range("N138").Value = -65235524
# Row 141
$ws.Range("H141").Value = 3740.7144
$ws.Range("I141").Value = 3644.2307
$ws.Range("K141").Value = 10932.6921
$ws.Range("M141").Value = -5752.6921

# ---- Sheet 2: ARM ----
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 4686.4287
$ws.Range("I32").Value = 4655.827
$ws.Range("K32").Value = 4655.827
$ws.Range("M32").Value = -4368.827
# Row 45
$ws.Range("H45").Value = 116587.5
$ws.Range("I45").Value = 171131.5
$ws.Range("K45").Value = 171131.5
$ws.Range("M45").Value = -170754.5
# Row 61
$ws.Range("H61").Value = 13895190
$ws.Range("I61").Value = 15157025
$ws.Range("K61").Value = 15157025
$ws.Range("M61").Value = -15156813
# Row 74
$ws.Range("H74").Value = 2434.5
$ws.Range("J74").Value = 5166.6665
$ws.Range("L74").Value = 5166.6665
$ws.Range("N74").Value = -6914.6665
# Row 77
$ws.Range("H77").Value = 2434.5
$ws.Range("J77").Value = 5166.6665
$ws.Range("L77").Value = 25833.3325
$ws.Range("N77").Value = -34569.3325
# Row 110
$ws.Range("H110").Value = 8857.333000000001
$ws.Range("I110").Value = 7411.125
$ws.Range("J110").Value = 11749.75
$ws.Range("K110").Value = 7411.125
$ws.Range("L110").Value = 11749.75
$ws.Range("M110").Value = -5366.125
$ws.Range("N110").Value = -15839.75
# Row 122
$ws.Range("H122").Value = 2552.3333
$ws.Range("I122").Value = 2511.6365
$ws.Range("K122").Value = 7534.9095
$ws.Range("M122").Value = -5084.9095
# Row 132
$ws.Range("H132").Value = 5892.787
$ws.Range("I132").Value = 5531.8535
$ws.Range("J132").Value = 8359.166999999999
$ws.Range("K132").Value = 16595.5605
$ws.Range("L132").Value = 25077.501
$ws.Range("M132").Value = -14065.5605
$ws.Range("N132").Value = -30137.501
# Row 136
$ws.Range("H136").Value = 13895190
$ws.Range("I136").Value = 15157025
$ws.Range("K136").Value = 45471075
$ws.Range("M136").Value = -45468525

# ---- Sheet 3: BSM ----
$ws = $wb.Worksheets.Item(3)
# Row 20
$ws.Range("H20").Value = 2777.6667
$ws.Range("I20").Value = 2486.5557
$ws.Range("J20").Value = 3651
$ws.Range("K20").Value = 2486.5557
$ws.Range("L20").Value = 3651
$ws.Range("M20").Value = -2239.5557
$ws.Range("N20").Value = -4145
# Row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# Row 86
$ws.Range("H86").Value = 2770.4285
$ws.Range("I86").Value = 2499.5
$ws.Range("K86").Value = 2499.5
$ws.Range("M86").Value = -1376.5
# Row 89
$ws.Range("H89").Value = 2770.4285
$ws.Range("I89").Value = 2499.5
$ws.Range("K89").Value = 12497.5
$ws.Range("M89").Value = -6881.5
# Row 94
$ws.Range("H94").Value = 1953.7667
$ws.Range("I94").Value = 1814.5454
$ws.Range("J94").Value = 2336.625
$ws.Range("K94").Value = 1814.5454
$ws.Range("L94").Value = 2336.625
$ws.Range("M94").Value = -1363.5454
$ws.Range("N94").Value = -3238.625
# Row 105
$ws.Range("H105").Value = 3467.0908
$ws.Range("I105").Value = 3155.625
$ws.Range("K105").Value = 3155.625
$ws.Range("M105").Value = -1408.625
# Row 107
$ws.Range("H107").Value = 5499.222
$ws.Range("I107").Value = 4332.5
$ws.Range("J107").Value = 7832.6665
$ws.Range("K107").Value = 4332.5
$ws.Range("L107").Value = 7832.6665
$ws.Range("M107").Value = -2412.5
$ws.Range("N107").Value = -11672.6665
# Row 132
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
# Row 134
$ws.Range("H134").Value = 2778.9666
$ws.Range("I134").Value = 2769.1428
$ws.Range("K134").Value = 8307.428400000001
$ws.Range("M134").Value = -5772.428400000001

# ---- Sheet 4: CRP ----
$ws = $wb.Worksheets.Item(4)
# Row 16
$ws.Range("H16").Value = 3561.4375
$ws.Range("I16").Value = 2805.8
$ws.Range("J16").Value = 4820.8335
$ws.Range("K16").Value = 2805.8
$ws.Range("L16").Value = 4820.8335
$ws.Range("M16").Value = -2518.8
$ws.Range("N16").Value = -5394.8335
# Row 19
$ws.Range("H19").Value = 508.69232
$ws.Range("I19").Value = 466.41666
$ws.Range("J19").Value = 1016
$ws.Range("K19").Value = 466.41666
$ws.Range("L19").Value = 1016
$ws.Range("M19").Value = -296.41666
$ws.Range("N19").Value = -1356
# Row 21
$ws.Range("H21").Value = 9999
$ws.Range("I21").Value = 9999
$ws.Range("K21").Value = 9999
$ws.Range("M21").Value = -9764
# Row 22
$ws.Range("H22").Value = 312.5
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = 50
# Row 24
$ws.Range("H24").Value = 508.69232
$ws.Range("I24").Value = 466.41666
$ws.Range("J24").Value = 1016
$ws.Range("K24").Value = 466.41666
$ws.Range("L24").Value = 1016
$ws.Range("M24").Value = -296.41666
$ws.Range("N24").Value = -1356
# Row 88
$ws.Range("H88").Value = 13866.6
$ws.Range("J88").Value = 13866.6
$ws.Range("L88").Value = 13866.6
$ws.Range("N88").Value = -14678.6
# Row 91
$ws.Range("H91").Value = 13866.6
$ws.Range("J91").Value = 13866.6
$ws.Range("L91").Value = 13866.6
$ws.Range("N91").Value = -16674.6
# Row 92
$ws.Range("H92").Value = 61500
$ws.Range("J92").Value = 61500
$ws.Range("L92").Value = 61500
$ws.Range("N92").Value = -66492
# Row 94
$ws.Range("H94").Value = 2664.6667
$ws.Range("J94").Value = 3999
$ws.Range("L94").Value = 3999
$ws.Range("N94").Value = -4901
# Row 107
$ws.Range("H107").Value = 481.95456
$ws.Range("I107").Value = 383.07693
$ws.Range("K107").Value = 383.07693
$ws.Range("M107").Value = 1536.92307
# Row 111
$ws.Range("H111").Value = 77919.5
$ws.Range("J111").Value = 77919.5
$ws.Range("L111").Value = 77919.5
$ws.Range("N111").Value = -86099.5
# Row 113
$ws.Range("H113").Value = 3561.4375
$ws.Range("I113").Value = 2805.8
$ws.Range("J113").Value = 4820.8335
$ws.Range("K113").Value = 2805.8
$ws.Range("L113").Value = 4820.8335
$ws.Range("M113").Value = -635.8000000000002
$ws.Range("N113").Value = -9160.833500000001
# Row 122
$ws.Range("H122").Value = 3678.4614
$ws.Range("I122").Value = 3485.7646
$ws.Range("K122").Value = 10457.2938
$ws.Range("M122").Value = -8007.293799999999
# Row 132
$ws.Range("H132").Value = 3570.04
$ws.Range("I132").Value = 2482.95
$ws.Range("J132").Value = 7918.4
$ws.Range("K132").Value = 7448.849999999999
$ws.Range("L132").Value = 23755.2
$ws.Range("M132").Value = -4918.849999999999
$ws.Range("N132").Value = -28815.2
# Row 134
$ws.Range("H134").Value = 5748.65
$ws.Range("I134").Value = 4779.25
$ws.Range("J134").Value = 9626.25
$ws.Range("K134").Value = 14337.75
$ws.Range("L134").Value = 28878.75
$ws.Range("M134").Value = -11802.75
$ws.Range("N134").Value = -33948.75

# ---- Sheet 5: CUL ----
$ws = $wb.Worksheets.Item(5)
# Row 3
$ws.Range("H3").Value = 7687.2
$ws.Range("I3").Value = 7421.75
$ws.Range("J3").Value = 7990.5713
$ws.Range("K3").Value = 22265.25
$ws.Range("L3").Value = 23971.7139
$ws.Range("M3").Value = -22153.25
$ws.Range("N3").Value = -24195.7139
# Row 8
$ws.Range("H8").Value = 561.2222
$ws.Range("I8").Value = 561.2222
$ws.Range("K8").Value = 1683.6666
$ws.Range("M8").Value = -1544.6666
# Row 17
$ws.Range("H17").Value = 52.166668
$ws.Range("I17").Value = 42.6
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 127.8
$ws.Range("L17").Value = 300
$ws.Range("M17").Value = 41.19999999999999
$ws.Range("N17").Value = -638
# Row 23
$ws.Range("H23").Value = 8570.583000000001
$ws.Range("I23").Value = 164
$ws.Range("J23").Value = 20339.8
$ws.Range("K23").Value = 492
$ws.Range("L23").Value = 61019.39999999999
$ws.Range("M23").Value = -257
$ws.Range("N23").Value = -61489.39999999999
# Row 46
$ws.Range("H46").Value = 65
$ws.Range("I46").Value = 65
$ws.Range("K46").Value = 195
$ws.Range("M46").Value = -104
# Row 68
$ws.Range("H68").Value = 71430100
$ws.Range("I68").Value = 83334984
$ws.Range("K68").Value = 250004952
$ws.Range("M68").Value = -250004141
# Row 71
$ws.Range("H71").Value = 71430100
$ws.Range("I71").Value = 83334984
$ws.Range("K71").Value = 750014856
$ws.Range("M71").Value = -750010800
# Row 92
$ws.Range("H92").Value = 165.45454
$ws.Range("J92").Value = 148.88889
$ws.Range("L92").Value = 446.66667
$ws.Range("N92").Value = -2942.66667
# Row 132
$ws.Range("H132").Value = 20834760
$ws.Range("I132").Value = 31251032
$ws.Range("J132").Value = 2218.375
$ws.Range("K132").Value = 281259288
$ws.Range("L132").Value = 19965.375
$ws.Range("M132").Value = -281256758
$ws.Range("N132").Value = -25025.375

# ---- Sheet 6: GSM ----
$ws = $wb.Worksheets.Item(6)
# Row 70
$ws.Range("H70").Value = 9574.75
$ws.Range("I70").Value = 7749.25
$ws.Range("J70").Value = 10487.5
$ws.Range("K70").Value = 7749.25
$ws.Range("L70").Value = 10487.5
$ws.Range("M70").Value = -7479.25
$ws.Range("N70").Value = -11027.5
# Row 73
$ws.Range("H73").Value = 9574.75
$ws.Range("I73").Value = 7749.25
$ws.Range("J73").Value = 10487.5
$ws.Range("K73").Value = 7749.25
$ws.Range("L73").Value = 10487.5
$ws.Range("M73").Value = -6813.25
$ws.Range("N73").Value = -12359.5
# Row 80
$ws.Range("H80").Value = 4539.75
$ws.Range("I80").Value = 4246.8
$ws.Range("J80").Value = 5028
$ws.Range("K80").Value = 4246.8
$ws.Range("L80").Value = 5028
$ws.Range("M80").Value = -3248.8
$ws.Range("N80").Value = -7024
# Row 83
$ws.Range("H83").Value = 4539.75
$ws.Range("I83").Value = 4246.8
$ws.Range("J83").Value = 5028
$ws.Range("K83").Value = 21234
$ws.Range("L83").Value = 25140
$ws.Range("M83").Value = -16242
$ws.Range("N83").Value = -35124
# Row 97
$ws.Range("H97").Value = 3524.4
$ws.Range("I97").Value = 975.0833
$ws.Range("J97").Value = 7348.375
$ws.Range("K97").Value = 975.0833
$ws.Range("L97").Value = 7348.375
$ws.Range("M97").Value = -479.0833
$ws.Range("N97").Value = -8340.375
# Row 102
$ws.Range("H102").Value = 7744.1665
$ws.Range("I102").Value = 6893
$ws.Range("K102").Value = 6893
$ws.Range("M102").Value = -5271
# Row 122
$ws.Range("H122").Value = 3047.2273
$ws.Range("I122").Value = 2741.5264
$ws.Range("J122").Value = 4983.3335
$ws.Range("K122").Value = 8224.5792
$ws.Range("L122").Value = 14950.0005
$ws.Range("M122").Value = -5774.5792
$ws.Range("N122").Value = -19850.0005
# Row 126
$ws.Range("H126").Value = 4453.6665
$ws.Range("I126").Value = 3927.182
$ws.Range("K126").Value = 11781.546
$ws.Range("M126").Value = -9311.545999999998
# Row 132
$ws.Range("H132").Value = 5600.227
$ws.Range("I132").Value = 4120.273
$ws.Range("K132").Value = 12360.819
$ws.Range("M132").Value = -9830.819

# ---- Sheet 7: LTW ----
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Range("H22").Value = 4149.8887
$ws.Range("I22").Value = 3058.3333
$ws.Range("J22").Value = 6333
$ws.Range("K22").Value = 3058.3333
$ws.Range("L22").Value = 6333
$ws.Range("M22").Value = -2763.3333
$ws.Range("N22").Value = -6923
# Row 27
$ws.Range("H27").Value = 4149.8887
$ws.Range("I27").Value = 3058.3333
$ws.Range("J27").Value = 6333
$ws.Range("K27").Value = 3058.3333
$ws.Range("L27").Value = 6333
$ws.Range("M27").Value = -2951.3333
$ws.Range("N27").Value = -6547
# Row 55
$ws.Range("H55").Value = 971.6667
$ws.Range("I55").Value = 1007.5
$ws.Range("K55").Value = 1007.5
$ws.Range("M55").Value = -834.5
# Row 68
$ws.Range("H68").Value = 17374.75
$ws.Range("I68").Value = 17333
$ws.Range("K68").Value = 17333
$ws.Range("M68").Value = -16584
# Row 71
$ws.Range("H71").Value = 17374.75
$ws.Range("I71").Value = 17333
$ws.Range("K71").Value = 86665
$ws.Range("M71").Value = -82921
# Row 82
$ws.Range("H82").Value = 8974.27
$ws.Range("I82").Value = 11880.077
$ws.Range("J82").Value = 6068.4614
$ws.Range("K82").Value = 11880.077
$ws.Range("L82").Value = 6068.4614
$ws.Range("M82").Value = -11519.077
$ws.Range("N82").Value = -6790.4614
# Row 85
$ws.Range("H85").Value = 8974.27
$ws.Range("I85").Value = 11880.077
$ws.Range("J85").Value = 6068.4614
$ws.Range("K85").Value = 11880.077
$ws.Range("L85").Value = 6068.4614
$ws.Range("M85").Value = -10632.077
$ws.Range("N85").Value = -8564.4614
# Row 98
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
# Row 100
$ws.Range("H100").Value = 2946069.8
$ws.Range("I100").Value = 4550335
$ws.Range("J100").Value = 4916.5
$ws.Range("K100").Value = 4550335
$ws.Range("L100").Value = 4916.5
$ws.Range("M100").Value = -4549794
$ws.Range("N100").Value = -5998.5
# Row 122
$ws.Range("H122").Value = 2622.3076
$ws.Range("I122").Value = 2634.818
$ws.Range("J122").Value = 2553.5
$ws.Range("K122").Value = 7904.454000000001
$ws.Range("L122").Value = 7660.5
$ws.Range("M122").Value = -5454.454000000001
$ws.Range("N122").Value = -12560.5
# Row 132
$ws.Range("H132").Value = 7465
$ws.Range("I132").Value = 7551.7617
$ws.Range("J132").Value = 7161.3335
$ws.Range("K132").Value = 22655.2851
$ws.Range("L132").Value = 21484.0005
$ws.Range("M132").Value = -20125.2851
$ws.Range("N132").Value = -26544.0005

# ---- Sheet 8: WVR ----
$ws = $wb.Worksheets.Item(8)
# Row 12
$ws.Range("H12").Value = 10999.667
$ws.Range("I12").Value = 11499.5
$ws.Range("K12").Value = 11499.5
$ws.Range("M12").Value = -11357.5
# Row 80
$ws.Range("H80").Value = 15000
$ws.Range("J80").Value = 15000
$ws.Range("L80").Value = 15000
$ws.Range("N80").Value = -16996
# Row 81
$ws.Range("H81").Value = 2590
$ws.Range("I81").Value = 1575.125
$ws.Range("J81").Value = 5296.3335
$ws.Range("K81").Value = 3150.25
$ws.Range("L81").Value = 10592.667
$ws.Range("M81").Value = -2089.25
$ws.Range("N81").Value = -12714.667
# Row 83
$ws.Range("H83").Value = 15000
$ws.Range("J83").Value = 15000
$ws.Range("L83").Value = 45000
$ws.Range("N83").Value = -54984
# Row 84
$ws.Range("H84").Value = 2590
$ws.Range("I84").Value = 1575.125
$ws.Range("J84").Value = 5296.3335
$ws.Range("K84").Value = 15751.25
$ws.Range("L84").Value = 52963.335
$ws.Range("M84").Value = -10447.25
$ws.Range("N84").Value = -63571.335
# Row 96
$ws.Range("H96").Value = 2404.923
$ws.Range("I96").Value = 1981.4
$ws.Range("K96").Value = 1981.4
$ws.Range("M96").Value = -608.4000000000001
# Row 113
$ws.Range("H113").Value = 584.4583
$ws.Range("I113").Value = 571.1579
$ws.Range("K113").Value = 1713.4737
$ws.Range("M113").Value = 456.5263
# Row 122
$ws.Range("H122").Value = 4002.5186
$ws.Range("I122").Value = 2291.8572
$ws.Range("J122").Value = 5844.769
$ws.Range("K122").Value = 6875.571599999999
$ws.Range("L122").Value = 17534.307
$ws.Range("M122").Value = -4425.571599999999
$ws.Range("N122").Value = -22434.307
# Row 136
$ws.Range("H136").Value = 4057.3
$ws.Range("I136").Value = 2889.9333
$ws.Range("K136").Value = 8669.7999
$ws.Range("M136").Value = -6119.7999

